# Applies the Constraints1.xlsx edit:
#  1. Reorders the existing person sheets (yoni, tair, asaf -> tair, asaf, yoni)
#     and adds five new person sheets (adir, stav, rotem, michal, emilia),
#     each built from the same weekly Morning/Evening "NO"-availability
#     template used by the existing person sheets.
#  2. Updates the "shifts" master roster sheet with the new assignment grid.

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: drop the three existing per-person sheets; they will be rebuilt in
# the correct order together with the five brand-new person sheets.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("yoni").Delete()
$wb.Worksheets.Item("tair").Delete()
$wb.Worksheets.Item("asaf").Delete()

# ---------------------------------------------------------------------------
# Step 2: (re)create the person sheets, in order, after the "shifts" sheet.
# Each sheet gets the same little template:
#   B1:H1 -> Sunday..Saturday
#   A2 -> Morning, A3 -> Evening
#   one "NO" mark in row 2 (morning) and one in row 3 (evening)
# ---------------------------------------------------------------------------
$days = @("Sunday","Monday","Tuesday","Wednesday","Thursday","Friday","Saturday")
$dayCols = @("B","C","D","E","F","G","H")

$people = @(
    @{ Name = "tair";   Morning = "E2"; Evening = "F3" },
    @{ Name = "asaf";   Morning = "B2"; Evening = "C3" },
    @{ Name = "yoni";   Morning = "D2"; Evening = "D3" },
    @{ Name = "adir";   Morning = "B2"; Evening = "C3" },
    @{ Name = "stav";   Morning = "G2"; Evening = "H3" },
    @{ Name = "rotem";  Morning = "G2"; Evening = "H3" },
    @{ Name = "michal"; Morning = "G2"; Evening = "H3" },
    @{ Name = "emilia"; Morning = "B2"; Evening = "C3" }
)

foreach ($person in $people) {
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ps = $wb.Worksheets.Add($null, $last)
    $ps.Name = $person.Name

    for ($i = 0; $i -lt $dayCols.Length; $i++) {
        $ps.Range($dayCols[$i] + "1").Value = $days[$i]
    }
    $ps.Range("A2").Value = "Morning"
    $ps.Range("A3").Value = "Evening"

    $ps.Range($person.Morning).Value = "NO"
    $ps.Range($person.Evening).Value = "NO"
}

# ---------------------------------------------------------------------------
# Step 3: update the "shifts" master roster sheet with the new grid of names.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("shifts")

$ws.Range("B2").Value = "rotem"
$ws.Range("C2").Value = "yoni"
$ws.Range("D2").Value = "stav"
$ws.Range("E2").Value = "adir"
$ws.Range("F2").Value = "adir"
$ws.Range("G2").Value = "tair"

$ws.Range("B3").Value = "yoni"
$ws.Range("C3").Value = "adir"
$ws.Range("D3").Value = "adir"
$ws.Range("E3").Value = "yoni"
$ws.Range("F3").Value = "stav"
$ws.Range("G3").Value = "adir"

$ws.Range("B4").Value = "michal"
$ws.Range("C4").Value = "emilia"
$ws.Range("D4").Value = "michal"
$ws.Range("E4").Value = "emilia"
$ws.Range("F4").Value = "michal"
$ws.Range("G4").Value = "emilia"

$ws.Range("B5").Value = "stav"
$ws.Range("C5").Value = "tair"
$ws.Range("D5").Value = "adir"
$ws.Range("E5").Value = "rotem"
$ws.Range("F5").Value = "adir"
$ws.Range("H5").Value = "yoni"

$ws.Range("B6").Value = "yoni"
$ws.Range("C6").Value = "stav"
$ws.Range("D6").Value = "tair"
$ws.Range("E6").Value = "tair"
$ws.Range("F6").Value = "asaf"
$ws.Range("H6").Value = "tair"

$ws.Range("B7").Value = "emilia"
$ws.Range("C7").Value = "michal"
$ws.Range("D7").Value = "emilia"
$ws.Range("E7").Value = "michal"
$ws.Range("F7").Value = "emilia"
$ws.Range("H7").Value = "michal"

$ws.Select()
$ws.Range("A1").Select()
